# Add 3 new groups of data (cxq, hyy, hzj) as new columns E, F, G.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers. Entry order matches the shared-string table order in the target
# workbook (cxq, then hyy, then hzj) even though the columns end up as
# E=cxq, F=hzj, G=hyy.
$ws.Range("E1").Value = "cxq6hz_20170224_144343_ASIC_EEG"
$ws.Range("G1").Value = "hyy-调节6Hz_20170306_110203_ASIC_EEG"
$ws.Range("F1").Value = "hzj-调节6Hz_20170220_113105_ASIC_EEG"

# Row 2 values
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1

# Row 3 values
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.98709677419354835
$ws.Range("G3").Value = 0.98634812286689422

# Leave the selection on column F (whole column), matching the author's
# final UI state after entering the data.
$ws.Range("F1:F1048576").Select()
